$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings (e.g. "1.001", "0.3300", "15.10") keep their exact original
# text representation instead of being coerced into floating point
# numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.339.71"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").Value = "1.936.28"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "250.92"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "0.7099"
$ws.Range("E6").Value = "  -4.21%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "0.3300"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "27.58"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "0.07311"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").Value = "0.8059"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").Value = "0.08071"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "1.935.44"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "5.484"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "94.61"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "15.10"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "30.341.23"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").Value = "253.26"
$ws.Range("E18").Value = "  -4.80%  "
$ws.Range("D19").Value = "0.000008208"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "5.812"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "2.190.27"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "9.746"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "164.08"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.30"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.340"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "0.1288"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("D30").Value = "1.349"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "1.541"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "4.419"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "4.166"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").Value = "1.263"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").Value = "0.7472"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("D37").Value = "2.773"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "0.01971"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "2.807"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("D40").Value = "78.85"
$ws.Range("E40").Value = "  -6.51%  "
$ws.Range("D41").Value = "6.424"
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("D42").Value = "0.4528"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").Value = "2.014"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").Value = "0.8467"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "101.76"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").Value = "9.733"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").Value = "7.438"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("D49").Value = "36.68"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "0.4178"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "0.06032"
$ws.Range("E51").Value = "  -0.78%  "

# Restore the (unformatted) default look of column D now that the
# text values have been safely written, to match the original file
# which had no explicit number format on these cells.
$ws.Range("D2:D51").ClearFormats()
